$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2,  "trainingaudio/25_tapapi1.wav", "pngimages/25_apple.png"),
    @(3,  "trainingaudio/08_tipako2.wav", "pngimages/08_bell.png"),
    @(4,  "trainingaudio/09_tipata2.wav", "pngimages/09_plane.png"),
    @(5,  "trainingaudio/13_kopopi1.wav", "pngimages/13_toast.png"),
    @(6,  "trainingaudio/03_kikita3.wav", "pngimages/03_box.png"),
    @(7,  "trainingaudio/10_tokiti1.wav", "pngimages/10_backpack.png"),
    @(8,  "trainingaudio/11_tokiko1.wav", "pngimages/11_compass.png"),
    @(9,  "trainingaudio/24_takopa1.wav", "pngimages/24_banana.png"),
    @(10, "trainingaudio/17_kotako2.wav", "pngimages/17_cracker.png")
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
